$wb = $excel.ActiveWorkbook

# --- OFF sheet: Road (row 3) updates ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 396
$wsOff.Range("C3").Value = 270
$wsOff.Range("D3").Value = 67
$wsOff.Range("F3").Value = 8

# --- DEF sheet: Road (row 3) updates ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 524
$wsDef.Range("C3").Value = 374
$wsDef.Range("D3").Value = 110
$wsDef.Range("G3").Value = 8
